# Generate Report for Handback
# The f3e2fe82-073d-4902-8608-044de7c6793f.md file has now been handed back
# (it is in sync with en-US), so update the status / timestamps / error
# details that the localization-status report tracks for that row across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the f3e2fe82-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet: row 3 is the f3e2fe82-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-09-02 18:54:37"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet: row 3 is the f3e2fe82-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-09-02 18:54:45"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
